# Update NATMI Ifnk-Ifnar2 LR-pair sheet with re-computed TPM-based stats.
# New data adds "ECs" as its own sending cluster (rows 2-4) in addition to
# the existing "MuSCs" sending cluster (rows 5-7, moved down and recalculated).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Ifnk", "Ifnar2", "ECs", 2, 0.6666666666666666, 0.301559, 0.904677, 0.2570607550566701, 0.2570607550566701, 3, 1, 16.69780766666667, 50.093423, 0.3087131902856588, 0.3087131902856588, 5.035374182152332, 45.318367639371, 0.07935804579078493, 0.07935804579078494),
    @("ECs", "Ifnk", "Ifnar2", "FAPs", 2, 0.6666666666666666, 0.301559, 0.904677, 0.2570607550566701, 0.2570607550566701, 3, 1, 27.876023, 83.628069, 0.5153787949052554, 0.5153787949052554, 8.406265619856999, 75.65639057871299, 0.1324836621585417, 0.1324836621585417),
    @("ECs", "Ifnk", "Ifnar2", "MuSCs", 2, 0.6666666666666666, 0.301559, 0.904677, 0.2570607550566701, 0.2570607550566701, 3, 1, 9.514586000000001, 28.543758, 0.1759080148090858, 0.1759080148090857, 2.869209039574, 25.822881356166, 0.04521904710734349, 0.04521904710734349),
    @("MuSCs", "Ifnk", "Ifnar2", "ECs", 3, 1, 0.8715450000000001, 2.614635, 0.7429392449433299, 0.7429392449433299, 3, 1, 16.69780766666667, 50.093423, 0.3087131902856588, 0.3087131902856588, 14.552890782845, 130.976017045605, 0.2293551444948739, 0.2293551444948739),
    @("MuSCs", "Ifnk", "Ifnar2", "FAPs", 3, 1, 0.8715450000000001, 2.614635, 0.7429392449433299, 0.7429392449433299, 3, 1, 27.876023, 83.628069, 0.5153787949052554, 0.5153787949052554, 24.295208465535, 218.656876189815, 0.3828951327467137, 0.3828951327467137),
    @("MuSCs", "Ifnk", "Ifnar2", "MuSCs", 3, 1, 0.8715450000000001, 2.614635, 0.7429392449433299, 0.7429392449433299, 3, 1, 9.514586000000001, 28.543758, 0.1759080148090858, 0.1759080148090857, 8.292389855370002, 74.63150869833002, 0.1306889677017423, 0.1306889677017422)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
